# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (E) / "Valor Mora" (F) table on Hoja1 rows 16-55 is
# refreshed with the next block of periods: the period codes roll forward
# by one month (1607..1910 instead of 1607..1812 shifted) and the mora
# value attached to each period updates to match the new period ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Periodo Mora" labels (column E) for rows 16-55, in order.
$periodos = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910"
)

$firstRow = 16

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]

    if ($row -le 41) {
        $ws.Cells.Item($row, 6).Value = 25774
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }
}
